# Traceability matrix refresh (R0.WP2 close-out: T-41 + T-44 complete,
# work package now 27/27 points - 100%).
#
# - "Traceability Matrix" sheet: Status column (G) updated per task, and
#   Last Updated column (J) bumped to the new regeneration date for every
#   data row.
# - "Summary" sheet: Last Generated timestamp refreshed to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Traceability Matrix")

$statuses = [ordered]@{
    2  = "Planificado"
    3  = "Planificado"
    4  = "Completado"
    5  = "Completado"
    6  = "Completado"
    7  = "En Progreso"
    8  = "Completado"
    9  = "Completado"
    10 = "Planificado"
    11 = "Planificado"
    12 = "Completado"
    13 = "Planificado"
}

$lastUpdated = "2025-07-01"

# Keep the "Last Updated" column as literal text (it was authored as text,
# not an Excel date serial) by formatting it before writing the value.
$ws.Range("J2:J13").NumberFormat = "@"

foreach ($row in $statuses.Keys) {
    $ws.Range("G$row").Value = $statuses[$row]
    $ws.Range("J$row").Value = $lastUpdated
}

$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = "2025-07-01T15:33:51.951Z"
